$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "N/A" values in column F for rows 4-8 (matching style of column E)
$ws.Range("F4").Value = "N/A"
$ws.Range("F5").Value = "N/A"
$ws.Range("F6").Value = "N/A"
$ws.Range("F7").Value = "N/A"
$ws.Range("F8").Value = "N/A"

# Copy style from column E to column F for rows 4-8
$ws.Range("E4:E8").Copy()
$ws.Range("F4:F8").PasteSpecial(-4122)

# Update the selected cell/range
$ws.Range("B8").Select()
